# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders / refreshes the "Estado de Cuenta" detail table (rows 16-41) so
# that, for every mora period (2004..2012, 2101..2104, ascending), the two
# workers (ELMIS ISABEL SUAREZ ROJAS / ZUNILDA MARIA BOLIVAR LEON) appear
# back to back, and corrects the "Valor Mora" (column F) so the 30430
# figure belongs to period 2104 (it had erroneously been on period 2104's
# first occurrence before, i.e. rows 16 and 29) instead of always being on
# the very first row of each worker's block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104")

$workers = @(
    @("45454358", "ELMIS ISABEL SUAREZ ROJAS"),
    @("45476669", "ZUNILDA MARIA BOLIVAR LEON")
)

$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        if ($period -eq "2104") {
            $valorMora = 30430
        } else {
            $valorMora = 35112
        }

        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $worker[0]
        $ws.Cells.Item($row, 4).Value = $worker[1]
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = $valorMora
        $ws.Cells.Item($row, 7).Value = 877803

        $row = $row + 1
    }
}
